# The lab moved the source images from a Mac-style folder ("Extinction/…")
# to a Windows folder the author copied in from their own machine
# ("Extinction\…"). Update the two distinct path strings used in column A
# (rows 2-11 use the CS+ image, rows 12-21 use the CS- image) so the
# backslash-style Windows path is used throughout, matching the "Added
# original huang experiment folder" re-pointing of the source folder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPlus  = "Extinction/CS+2.BMP"
$newPlus  = "Extinction\CS+2.BMP"
$oldMinus = "Extinction/CS-2.BMP"
$newMinus = "Extinction\CS-2.BMP"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()
    if ($v -eq $oldPlus) {
        $cell.Value = $newPlus
    } elseif ($v -eq $oldMinus) {
        $cell.Value = $newMinus
    }
}
